# Auto commit 28-05-2025 11:31:48.36
#
# Adds two new template sheets ("Submersible Pump Installation" and
# "Compressor Pump Installation") at the end of the workbook. Each is a
# truncated copy of an existing template:
#   - "Submersible Pump Installation"  <- first 12 rows of "MWSS with Submersible Pump"
#   - "Compressor Pump Installation"   <- first 16 rows of "MWSS with Compressor Pump"
# Also tidies up a couple of stale sheet-view selections left over on the
# "150 mm TWC" / "200 mm TWC" tabs, and the row-1 height on the
# "MWSS with Compressor Pump" sheet.

$wb = $excel.ActiveWorkbook

$wsSubmersible = $wb.Worksheets.Item(1)   # "MWSS with Submersible Pump"
$wsCompressor  = $wb.Worksheets.Item(2)   # "MWSS with Compressor Pump"
$wsTwc150      = $wb.Worksheets.Item(5)   # "150 mm TWC"
$wsTwc200      = $wb.Worksheets.Item(6)   # "200 mm TWC"

# --- cosmetic fix on "MWSS with Compressor Pump" header row -----------------
$wsCompressor.Rows.Item(1).RowHeight = 29

# --- tidy the stray selections on the two TWC sheets -------------------------
$wsTwc150.Activate()
$wsTwc150.Range("A1:B1").Select()

$wsTwc200.Activate()
$wsTwc200.Range("A1:B1").Select()

# --- new sheet: Submersible Pump Installation --------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsNewSubmersible = $wb.Worksheets.Add($null, $lastSheet)
$wsNewSubmersible.Name = "Submersible Pump Installation"
$wsNewSubmersible.Range("A1:B12").Value = $wsSubmersible.Range("A1:B12").Value2
$wsNewSubmersible.Activate()
$wsNewSubmersible.Range("A1:B1").Select()

# --- new sheet: Compressor Pump Installation ----------------------------------
$wsNewCompressor = $wb.Worksheets.Add($null, $wsNewSubmersible)
$wsNewCompressor.Name = "Compressor Pump Installation"
$wsNewCompressor.Range("A1:B16").Value = $wsCompressor.Range("A1:B16").Value2
$wsNewCompressor.Activate()
$wsNewCompressor.Range("E5").Select()
